# Updates cryptos list data (columns B/C/D/E) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.963.68"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").Value = "'2.046.10"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").Value = "'251.54"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "'58.69"
$ws.Range("E7").Value = "  +3.89%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "'61.29"
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("D10").Value = "'0.387"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").Value = "'0.0787"
$ws.Range("E11").Value = "  +4.25%  "
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").Value = "'16.30"
$ws.Range("E13").Value = "  +5.66%  "
$ws.Range("D14").Value = "'2.340.91"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "'0.809"
$ws.Range("E15").Value = "  -5.52%  "
$ws.Range("E16").Value = "  +6.71%  "
$ws.Range("D17").Value = "'2.052.45"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").Value = "'36.895.01"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").Value = "'16.79"
$ws.Range("E19").Value = "  +14.94%  "
$ws.Range("D20").Value = "'75.00"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").Value = "'0.0₃0907"
$ws.Range("E21").Value = "  +6.11%  "
$ws.Range("D22").Value = "'5.44"
$ws.Range("E22").Value = "  +3.11%  "
$ws.Range("D23").Value = "'237.23"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("D26").Value = "'2.29"
$ws.Range("E26").Value = "  +13.07%  "
$ws.Range("D27").Value = "'169.16"
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("D28").Value = "'9.30"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Value = "'20.24"
$ws.Range("E29").Value = "  -3.06%  "
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  +5.63%  "
$ws.Range("D32").Value = "'4.75"
$ws.Range("E32").Value = "  +4.05%  "
$ws.Range("D33").Value = "'0.0621"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("D34").Value = "'4.48"
$ws.Range("E34").Value = "  +3.77%  "
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").Value = "'0.0874"
$ws.Range("E36").Value = "  -3.50%  "
$ws.Range("D37").Value = "'2.23"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  -4.06%  "
$ws.Range("E39").Value = "  +13.26%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "'17.91"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("D43").Value = "'1.14"
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("D44").Value = "'97.18"
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("E46").Value = "  +14.79%  "
$ws.Range("E47").Value = "  +5.63%  "
$ws.Range("D48").Value = "'1.286.38"
$ws.Range("E48").Value = "  -3.23%  "
$ws.Range("D49").Value = "'2.89"
$ws.Range("E49").Value = "  -2.01%  "
$ws.Range("D50").Value = "'6.77"
$ws.Range("E50").Value = "  -3.80%  "
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").Value = "'3.65"
$ws.Range("E51").Value = "  -15.97%  "
